$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Siglec1"
$ws.Range("C2").Value = "Spn"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.120277
$ws.Range("H2").Value = 0.360831
$ws.Range("I2").Value = 0.001062914843064901
$ws.Range("J2").Value = 0.00106414053856905
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.787414
$ws.Range("N2").Value = 11.362242
$ws.Range("O2").Value = 0.5877125485801681
$ws.Range("P2").Value = 0.587712548580168
$ws.Range("Q2").Value = 0.455538793678
$ws.Range("R2").Value = 4.099849143102
$ws.Range("S2").Value = 0.0006246883913413623
$ws.Range("T2").Value = 0.0006254087479698888

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Siglec1"
$ws.Range("C3").Value = "Spn"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.120277
$ws.Range("H3").Value = 0.360831
$ws.Range("I3").Value = 0.001062914843064901
$ws.Range("J3").Value = 0.00106414053856905
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.656916666666667
$ws.Range("N3").Value = 7.97075
$ws.Range("O3").Value = 0.412287451419832
$ws.Range("P3").Value = 0.4122874514198319
$ws.Range("Q3").Value = 0.3195659659166667
$ws.Range("R3").Value = 2.87609369325
$ws.Range("S3").Value = 0.0004382264517235387
$ws.Range("T3").Value = 0.0004387317905991608

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Siglec1"
$ws.Range("C4").Value = "Spn"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4782236666666666
$ws.Range("H4").Value = 1.434671
$ws.Range("I4").Value = 0.004226169871254865
$ws.Range("J4").Value = 0.004231043260167216
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.787414
$ws.Range("N4").Value = 11.362242
$ws.Range("O4").Value = 0.5877125485801681
$ws.Range("P4").Value = 0.587712548580168
$ws.Range("Q4").Value = 1.811231010264666
$ws.Range("R4").Value = 16.301079092382
$ws.Range("S4").Value = 0.002483773065767917
$ws.Range("T4").Value = 0.002486637217585817

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Siglec1"
$ws.Range("C5").Value = "Spn"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4782236666666666
$ws.Range("H5").Value = 1.434671
$ws.Range("I5").Value = 0.004226169871254865
$ws.Range("J5").Value = 0.004231043260167216
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.656916666666667
$ws.Range("N5").Value = 7.97075
$ws.Range("O5").Value = 0.412287451419832
$ws.Range("P5").Value = 0.4122874514198319
$ws.Range("Q5").Value = 1.270600430361111
$ws.Range("R5").Value = 11.43540387325
$ws.Range("S5").Value = 0.001742396805486948
$ws.Range("T5").Value = 0.001744406042581398

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Siglec1"
$ws.Range("C6").Value = "Spn"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 69.61810866666667
$ws.Range("H6").Value = 208.854326
$ws.Range("I6").Value = 0.6152308508518273
$ws.Range("J6").Value = 0.615940301559777
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.787414
$ws.Range("N6").Value = 11.362242
$ws.Range("O6").Value = 0.5877125485801681
$ws.Range("P6").Value = 0.587712548580168
$ws.Range("Q6").Value = 263.6725994176547
$ws.Range("R6").Value = 2373.053394758892
$ws.Range("S6").Value = 0.3615788913192727
$ws.Range("T6").Value = 0.3619958444029338

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Siglec1"
$ws.Range("C7").Value = "Spn"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 69.61810866666667
$ws.Range("H7").Value = 208.854326
$ws.Range("I7").Value = 0.6152308508518273
$ws.Range("J7").Value = 0.615940301559777
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.656916666666667
$ws.Range("N7").Value = 7.97075
$ws.Range("O7").Value = 0.412287451419832
$ws.Range("P7").Value = 0.4122874514198319
$ws.Range("Q7").Value = 184.9695132182778
$ws.Range("R7").Value = 1664.7256189645
$ws.Range("S7").Value = 0.2536519595325546
$ws.Range("T7").Value = 0.2539444571568432

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Siglec1"
$ws.Range("C8").Value = "Spn"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 42.55008066666667
$ws.Range("H8").Value = 127.650242
$ws.Range("I8").Value = 0.376024612471286
$ws.Range("J8").Value = 0.376458223573777
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.787414
$ws.Range("N8").Value = 11.362242
$ws.Range("O8").Value = 0.5877125485801681
$ws.Range("P8").Value = 0.587712548580168
$ws.Range("Q8").Value = 161.1547712180627
$ws.Range("R8").Value = 1450.392940962564
$ws.Range("S8").Value = 0.2209943833243696
$ws.Range("T8").Value = 0.2212492220105071

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Siglec1"
$ws.Range("C9").Value = "Spn"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 42.55008066666667
$ws.Range("H9").Value = 127.650242
$ws.Range("I9").Value = 0.376024612471286
$ws.Range("J9").Value = 0.376458223573777
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.656916666666667
$ws.Range("N9").Value = 7.97075
$ws.Range("O9").Value = 0.412287451419832
$ws.Range("P9").Value = 0.4122874514198319
$ws.Range("Q9").Value = 113.0520184912778
$ws.Range("R9").Value = 1017.4681664215
$ws.Range("S9").Value = 0.1550302291469165
$ws.Range("T9").Value = 0.1552090015632698

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Siglec1"
$ws.Range("C10").Value = "Spn"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.391011
$ws.Range("H10").Value = 0.782022
$ws.Range("I10").Value = 0.003455451962566825
$ws.Range("J10").Value = 0.00230629106770994
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.787414
$ws.Range("N10").Value = 11.362242
$ws.Range("O10").Value = 0.5877125485801681
$ws.Range("P10").Value = 0.587712548580168
$ws.Range("Q10").Value = 1.480920535554
$ws.Range("R10").Value = 8.885523213324
$ws.Range("S10").Value = 0.002030812479416492
$ws.Range("T10").Value = 0.001355436201171486

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Siglec1"
$ws.Range("C11").Value = "Spn"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.391011
$ws.Range("H11").Value = 0.782022
$ws.Range("I11").Value = 0.003455451962566825
$ws.Range("J11").Value = 0.00230629106770994
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.656916666666667
$ws.Range("N11").Value = 7.97075
$ws.Range("O11").Value = 0.412287451419832
$ws.Range("P11").Value = 0.4122874514198319
$ws.Range("Q11").Value = 1.03888364275
$ws.Range("R11").Value = 6.2333018565
$ws.Range("S11").Value = 0.001424639483150333
$ws.Range("T11").Value = 0.0009508548665384542
